$d = $word.ActiveDocument

# --- Edit 1: "Turkish" -> "Turkish (Baseline and Follow Up)" in the v0.1.5 list ---
$r1 = $d.Content
$r1.Find.Execute("Turkish") | Out-Null
$r1.Collapse(0)
$r1.InsertAfter(" (Baseline and Follow Up)")
$r1.Font.Name = "Arial"
$r1.Font.Size = 11
$r1.Font.NameBi = "Arial"

# --- Edit 2: add a new list item "Turkish (Current Forms)" after "Greek" ---
$r2 = $d.Content
$r2.Find.Execute("Greek") | Out-Null
$para = $r2.Paragraphs(1)
$startOfGreekPara = $para.Range.Start

$endRange = $para.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

# locate the newly-inserted (now empty) paragraph right after the "Greek" paragraph
$newPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs($i)
  if ($p.Range.Start -eq $startOfGreekPara) {
    $newPara = $d.Paragraphs($i + 1)
    break
  }
}

$newRange = $newPara.Range
$newRange.End = $newRange.End - 1
$newRange.Text = "Turkish (Current Forms)"
